$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B4/B5 activity values are swapped: B4 becomes "Programmazione", B5 becomes "Scrittura"
$ws.Range("B4").Value = "Programmazione"
$ws.Range("B5").Value = "Scrittura"

# New D4/D5 hour entries - columns D/E already carry the correct default
# column style (h:mm number format / borders), so plain value writes pick
# it up automatically.
# New E4/E5 note entries: E5's text is registered first in the shared
# strings table, so write it before E4 to match.
$ws.Range("E5").Value = "Scritti testi primi quattro funghi"
$ws.Range("D5").Value = 0.083333333333333329

$ws.Range("E4").Value = "Creata struttura gestione funghi etc"
$ws.Range("D4").Value = 0.10416666666666667

# Update selection to E4 (last active cell)
$ws.Range("E4").Select()
